$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 9261986
$ws.Range("I33").Value = 15626089
$ws.Range("K33").Value = 15626089
$ws.Range("M33").Value = -15625860
$ws.Range("H40").Value = 1243.75
$ws.Range("I40").Value = 987.5
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 987.5
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -812.5
$ws.Range("N40").Value = -1850
$ws.Range("H106").Value = 6368.6
$ws.Range("I106").Value = 6368.6
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 6368.6
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -5737.6
$ws.Range("H125").Value = 3500769.8
$ws.Range("I125").Value = 9098107
$ws.Range("K125").Value = 81882963
$ws.Range("M125").Value = -81880503
$ws.Range("H131").Value = 558327.3
$ws.Range("I131").Value = 716422.3
$ws.Range("K131").Value = 2149266.9
$ws.Range("M131").Value = -2144226.9
$ws.Range("H132").Value = 31135.5
$ws.Range("I132").Value = 32599
$ws.Range("K132").Value = 97797
$ws.Range("M132").Value = -95267
$ws.Range("H138").Value = 2237.6558
$ws.Range("J138").Value = 2073
$ws.Range("L138").Value = 6219
$ws.Range("N138").Value = -16499

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12418.608
$ws.Range("I74").Value = 1536.2222
$ws.Range("J74").Value = 51595.2
$ws.Range("K74").Value = 1536.2222
$ws.Range("L74").Value = 51595.2
$ws.Range("M74").Value = -662.2221999999999
$ws.Range("N74").Value = -53343.2
$ws.Range("H77").Value = 12418.608
$ws.Range("I77").Value = 1536.2222
$ws.Range("J77").Value = 51595.2
$ws.Range("K77").Value = 7681.111
$ws.Range("L77").Value = 257976
$ws.Range("M77").Value = -3313.111
$ws.Range("N77").Value = -266712
$ws.Range("H97").Value = 877.44446
$ws.Range("I97").Value = 877.44446
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 877.44446
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -381.44446

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7595117
$ws.Range("I20").Value = 14497848
$ws.Range("J20").Value = 34982.145
$ws.Range("K20").Value = 14497848
$ws.Range("L20").Value = 34982.145
$ws.Range("M20").Value = -14497601
$ws.Range("N20").Value = -35476.145
$ws.Range("H22").Value = 23811340
$ws.Range("I22").Value = 23811340
$ws.Range("K22").Value = 23811340
$ws.Range("M22").Value = -23811167
$ws.Range("H54").Value = 1083
$ws.Range("I54").Value = 1083
$ws.Range("K54").Value = 1083
$ws.Range("M54").Value = -599
$ws.Range("H115").Value = 80000
$ws.Range("J115").Value = 80000
$ws.Range("L115").Value = 80000
$ws.Range("N115").Value = -83134
$ws.Range("H134").Value = 32550.465
$ws.Range("I134").Value = 42221
$ws.Range("K134").Value = 126663
$ws.Range("M134").Value = -124128
$ws.Range("H137").Value = 88000
$ws.Range("J137").Value = 88000
$ws.Range("L137").Value = 88000
$ws.Range("N137").Value = -98200

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 14426.77
$ws.Range("I58").Value = 5868.6523
$ws.Range("K58").Value = 5868.6523
$ws.Range("M58").Value = -5665.6523
$ws.Range("H99").Value = 12605375
$ws.Range("I99").Value = 8905562
$ws.Range("J99").Value = 20005000
$ws.Range("K99").Value = 8905562
$ws.Range("L99").Value = 20005000
$ws.Range("M99").Value = -8904064
$ws.Range("N99").Value = -20007996
$ws.Range("H107").Value = 877.2727
$ws.Range("J107").Value = 1122.5
$ws.Range("L107").Value = 1122.5
$ws.Range("N107").Value = -4962.5
$ws.Range("H126").Value = 12605375
$ws.Range("I126").Value = 8905562
$ws.Range("J126").Value = 20005000
$ws.Range("K126").Value = 26716686
$ws.Range("L126").Value = 60015000
$ws.Range("M126").Value = -26714216
$ws.Range("N126").Value = -60019940
$ws.Range("H136").Value = 14426.77
$ws.Range("I136").Value = 5868.6523
$ws.Range("K136").Value = 17605.9569
$ws.Range("M136").Value = -15055.9569
$ws.Range("H141").Value = 174974.62
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 362.875
$ws.Range("I33").Value = 99
$ws.Range("K33").Value = 594
$ws.Range("M33").Value = -311
$ws.Range("H44").Value = 315.75
$ws.Range("I44").Value = 315.75
$ws.Range("K44").Value = 947.25
$ws.Range("M44").Value = -549.25
$ws.Range("H86").Value = 499.4138
$ws.Range("J86").Value = 409
$ws.Range("L86").Value = 1227
$ws.Range("N86").Value = -3599
$ws.Range("H89").Value = 499.4138
$ws.Range("J89").Value = 409
$ws.Range("L89").Value = 3681
$ws.Range("N89").Value = -15537
$ws.Range("H92").Value = 258.27274
$ws.Range("J92").Value = 299.8
$ws.Range("L92").Value = 899.4000000000001
$ws.Range("N92").Value = -3395.4

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 26864.072
$ws.Range("I20").Value = 8050
$ws.Range("J20").Value = 29999.75
$ws.Range("K20").Value = 8050
$ws.Range("L20").Value = 29999.75
$ws.Range("M20").Value = -7805
$ws.Range("N20").Value = -30489.75
$ws.Range("H23").Value = 6297.5
$ws.Range("I23").Value = 234.5
$ws.Range("J23").Value = 9329
$ws.Range("K23").Value = 234.5
$ws.Range("L23").Value = 9329
$ws.Range("M23").Value = -11.5
$ws.Range("N23").Value = -9775
$ws.Range("H24").Value = 27770.77
$ws.Range("J24").Value = 27770.77
$ws.Range("L24").Value = 27770.77
$ws.Range("N24").Value = -28116.77
$ws.Range("H102").Value = 6759781.5
$ws.Range("I102").Value = 7115480.5
$ws.Range("K102").Value = 7115480.5
$ws.Range("M102").Value = -7113858.5
$ws.Range("H126").Value = 3359840.8
$ws.Range("I126").Value = 1637418.1
$ws.Range("K126").Value = 4912254.300000001
$ws.Range("M126").Value = -4909784.300000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1541677.1
$ws.Range("J7").Value = 5972
$ws.Range("L7").Value = 5972
$ws.Range("N7").Value = -6196
$ws.Range("H32").Value = 9766.666999999999
$ws.Range("I32").Value = 2300
$ws.Range("K32").Value = 2300
$ws.Range("M32").Value = -1983
$ws.Range("H55").Value = 1830.8182
$ws.Range("I55").Value = 1312.4286
$ws.Range("K55").Value = 1312.4286
$ws.Range("M55").Value = -1139.4286
$ws.Range("H126").Value = 1541677.1
$ws.Range("J126").Value = 5972
$ws.Range("L126").Value = 17916
$ws.Range("N126").Value = -22856
$ws.Range("H132").Value = 1349671.9
$ws.Range("I132").Value = 2682.2173
$ws.Range("K132").Value = 8046.651899999999
$ws.Range("M132").Value = -5516.651899999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 47894.668
$ws.Range("J33").Value = 47894.668
$ws.Range("L33").Value = 47894.668
$ws.Range("N33").Value = -48394.668
$ws.Range("H36").Value = 47894.668
$ws.Range("J36").Value = 47894.668
$ws.Range("L36").Value = 47894.668
$ws.Range("N36").Value = -48394.668
$ws.Range("H64").Value = 560114
$ws.Range("J64").Value = 560114
$ws.Range("L64").Value = 560114
$ws.Range("N64").Value = -560610
$ws.Range("H67").Value = 560114
$ws.Range("J67").Value = 560114
$ws.Range("L67").Value = 560114
$ws.Range("N67").Value = -561830
$ws.Range("H81").Value = 1294.8572
$ws.Range("I81").Value = 1312.8
$ws.Range("K81").Value = 2625.6
$ws.Range("M81").Value = -1564.6
$ws.Range("H84").Value = 1294.8572
$ws.Range("I84").Value = 1312.8
$ws.Range("K84").Value = 13128
$ws.Range("M84").Value = -7824
$ws.Range("H96").Value = 1323.1666
$ws.Range("I96").Value = 1090
$ws.Range("J96").Value = 1689.5714
$ws.Range("K96").Value = 1090
$ws.Range("L96").Value = 1689.5714
$ws.Range("M96").Value = 283
$ws.Range("N96").Value = -4435.5714
$ws.Range("H126").Value = 6254322
$ws.Range("I126").Value = 5959.5454
$ws.Range("K126").Value = 17878.6362
$ws.Range("M126").Value = -15408.6362
$ws.Range("H132").Value = 298087.3
$ws.Range("I132").Value = 1730.4073
$ws.Range("K132").Value = 5191.2219
$ws.Range("M132").Value = -2661.2219
